$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 3 ("Fressen") for the new "Fahrrad" product.
# This shifts the old rows 3-9 down to 4-10 (mergeCells / summary rows follow).
$ws.Rows("3").Insert()

# The product rows use a banded (alternating) row style that is reapplied
# per final row position, not carried along with the row's original content.
# Row 2 ("Bett") keeps the "even" style-band; capture the "odd" style-band
# from row 4 (which currently still holds the pre-insert row-3 formatting)
# before we start redistributing it, then hand it out to every odd data row.
$ws.Range("A4:F4").Copy()
$ws.Range("A3:F3").PasteSpecial(-4122)
$ws.Rows("3").RowHeight = 20

# Even data rows (4, 6, 8) get the "even" style-band used by row 2.
$ws.Range("A2:F2").Copy()
$ws.Range("A4:F4").PasteSpecial(-4122)
$ws.Range("A6:F6").PasteSpecial(-4122)
$ws.Range("A8:F8").PasteSpecial(-4122)

# Odd data rows (5, 7) get the "odd" style-band now sitting on row 3.
$ws.Range("A3:F3").Copy()
$ws.Range("A5:F5").PasteSpecial(-4122)
$ws.Range("A7:F7").PasteSpecial(-4122)

# Fill in the new "Fahrrad" row's data.
$ws.Range("A3").Value = "Fahrrad"
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0.0
$ws.Range("E3").Value = "Katzenkiste"

# "Fressen" (now row 4) gets a "Katzentoilette" tag instead of "Ball".
$ws.Range("F4").Value = "Katzentoilette"

# "Toilette" (now row 8) gets its tag replaced with a new value.
$ws.Range("F8").Value = "sfasdfsdfadsafd"
